$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")

# CategoryCustomBasic
$ws.Cells.Item(42, 2).Value = "OptionBasicChance"
$ws.Cells.Item(43, 2).Value = "OptionBasicTimeScale"
$ws.Cells.Item(44, 2).Value = "OptionBasicDuration"
$ws.Cells.Item(45, 2).Value = "OptionBasicCooldown"
$ws.Cells.Item(46, 2).Value = "OptionBasicSmoothIn"
$ws.Cells.Item(47, 2).Value = "OptionBasicSmoothOut"
$ws.Cells.Item(48, 2).Value = "OptionBasicThirdPerson"

# CategoryCustomCritical
$ws.Cells.Item(51, 2).Value = "OptionCriticalChance"
$ws.Cells.Item(52, 2).Value = "OptionCriticalTimeScale"
$ws.Cells.Item(53, 2).Value = "OptionCriticalDuration"
$ws.Cells.Item(54, 2).Value = "OptionCriticalCooldown"
$ws.Cells.Item(55, 2).Value = "OptionCriticalSmoothIn"
$ws.Cells.Item(56, 2).Value = "OptionCriticalSmoothOut"
$ws.Cells.Item(57, 2).Value = "OptionCriticalThirdPerson"

# CategoryCustomDismemberment
$ws.Cells.Item(60, 2).Value = "OptionDismemberChance"
$ws.Cells.Item(61, 2).Value = "OptionDismemberTimeScale"
$ws.Cells.Item(62, 2).Value = "OptionDismemberDuration"
$ws.Cells.Item(63, 2).Value = "OptionDismemberCooldown"
$ws.Cells.Item(64, 2).Value = "OptionDismemberSmoothIn"
$ws.Cells.Item(65, 2).Value = "OptionDismemberSmoothOut"
$ws.Cells.Item(66, 2).Value = "OptionDismemberThirdPerson"

# CategoryCustomDecapitation
$ws.Cells.Item(69, 2).Value = "OptionDecapChance"
$ws.Cells.Item(70, 2).Value = "OptionDecapTimeScale"
$ws.Cells.Item(71, 2).Value = "OptionDecapDuration"
$ws.Cells.Item(72, 2).Value = "OptionDecapCooldown"
$ws.Cells.Item(73, 2).Value = "OptionDecapSmoothIn"
$ws.Cells.Item(74, 2).Value = "OptionDecapSmoothOut"
$ws.Cells.Item(75, 2).Value = "OptionDecapThirdPerson"

# CategoryCustomLastEnemy
$ws.Cells.Item(78, 2).Value = "OptionLastEnemyChance"
$ws.Cells.Item(79, 2).Value = "OptionLastEnemyTimeScale"
$ws.Cells.Item(80, 2).Value = "OptionLastEnemyDuration"
$ws.Cells.Item(81, 2).Value = "OptionLastEnemyCooldown"
$ws.Cells.Item(82, 2).Value = "OptionLastEnemySmoothIn"
$ws.Cells.Item(83, 2).Value = "OptionLastEnemySmoothOut"
$ws.Cells.Item(84, 2).Value = "OptionLastEnemyThirdPerson"

# CategoryCustomLastStand
$ws.Cells.Item(87, 2).Value = "OptionLastStandTimeScale"
$ws.Cells.Item(88, 2).Value = "OptionLastStandDuration"
$ws.Cells.Item(89, 2).Value = "OptionLastStandCooldown"
$ws.Cells.Item(90, 2).Value = "OptionLastStandSmoothIn"
$ws.Cells.Item(91, 2).Value = "OptionLastStandSmoothOut"

# CategoryCustomParry
$ws.Cells.Item(94, 2).Value = "OptionParryChance"
$ws.Cells.Item(95, 2).Value = "OptionParryTimeScale"
$ws.Cells.Item(96, 2).Value = "OptionParryDuration"
$ws.Cells.Item(97, 2).Value = "OptionParryCooldown"
$ws.Cells.Item(98, 2).Value = "OptionParrySmoothIn"
$ws.Cells.Item(99, 2).Value = "OptionParrySmoothOut"
